# Update "想去人数" (number of people interested) counts on the
# "展览" (Exhibition) and "全部类型" (All Types) worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 7491
$wsExhibit.Range("F5").Value = 236
$wsExhibit.Range("F6").Value = 1129
$wsExhibit.Range("F7").Value = 196
$wsExhibit.Range("F9").Value = 125
$wsExhibit.Range("F10").Value = 32

# --- Sheet "全部类型" (All Types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 7491
$wsAll.Range("F5").Value = 236
$wsAll.Range("F6").Value = 1129
$wsAll.Range("F7").Value = 196
$wsAll.Range("F10").Value = 125
$wsAll.Range("F11").Value = 32
